$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify that the "Total hours" figure in the labour-costs table is per person
$ws.Range("D12").Value = "Total hours/person"

# Column D needs to be a bit wider now that the header text is longer
# (target raw width ~21.71; 20.83 is the closest input that lands on the
# runtime's nearest achievable quantized column width)
$ws.Columns("D").ColumnWidth = 20.83

# Active selection moved to D13 as of this edit
$ws.Range("D13").Select()
